$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new Wins/Losses/Ties columns, matching the
# formatting of the existing header row (bold, bordered, centered) by
# copying the format from the neighboring header cell (AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for every data row (2-56) with the same
# win/loss/tie totals.
$ws.Range("AD2:AD56").Value = 78
$ws.Range("AE2:AE56").Value = 84
$ws.Range("AF2:AF56").Value = 0
